$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-25 04:07:19"
$wsZhCn.Range("G5").Value = "2016-01-25 04:08:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-25 04:07:30"
$wsDeDe.Range("G5").Value = "2016-01-25 04:08:23"
